# Auto price-data update: prepend a new day (2025-11-26) to the price
# history table and push the existing rows down by one, dropping the
# oldest visible date (2025-11-21) back in at the new bottom row with
# the same metal prices (prices are unchanged day over day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 2-6 down to rows 3-7, working bottom-up so data isn't
# clobbered before it has been copied. Column A holds date text like
# "2025-11-25"; force text formatting before writing so Excel doesn't
# reinterpret it as a date serial number.
for ($r = 6; $r -ge 2; $r--) {
    $dst = $r + 1

    $ws.Cells.Item($dst, 1).NumberFormat = "@"
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dst, 4).Value = $ws.Cells.Item($r, 4).Value2
}

# Write the newest day's data into row 2.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-11-26"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
